# TournRPG-239: add an "item discard" message/UI entry
#   - message sheet: new row for the "<val1>を捨てた" log message (green)
#   - ui sheet: new row for the "アイテム捨てる" menu label

$wb = $excel.ActiveWorkbook

# --- "message" sheet: append row 47 ---
$wsMsg = $wb.Worksheets.Item("message")

# Copy formatting (styles/borders/fill/number format) from the last existing
# data row so the new row matches the table's look.
$wsMsg.Range("A46:C46").Copy()
$wsMsg.Range("A47:C47").PasteSpecial(-4122)

$wsMsg.Cells.Item(47, 1).Formula = "=ROW()-2"
$wsMsg.Cells.Item(47, 2).Value = "<val1>を捨てた"
$wsMsg.Cells.Item(47, 3).Value = "green"
$wsMsg.Rows.Item(47).RowHeight = 20

# --- "ui" sheet: append row 20 ---
$wsUi = $wb.Worksheets.Item("ui")

$wsUi.Range("A19:B19").Copy()
$wsUi.Range("A20:B20").PasteSpecial(-4122)

$wsUi.Cells.Item(20, 1).Formula = "=ROW()-2"
$wsUi.Cells.Item(20, 2).Value = "アイテム捨てる"
$wsUi.Rows.Item(20).RowHeight = 20

Write-Host "Applied TournRPG-239 item-discard message/UI rows"
